# Apply updated ranking values to the "matrices" (F), "prolificid" (C),
# "name" (D), "index"/B, and "race" (G) columns on the active sheet.
# Rows are re-sorted/re-randomized and the matrices score is recomputed;
# a few rows swap their prolificid/name/index/race while the row's
# position (A = index, H = mat_rank) stays fixed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 13.37383182294894

# Row 3
$ws.Range("F3").Value = 13.16250246588412

# Row 4
$ws.Range("F4").Value = 8.396910801783761

# Row 5
$ws.Range("B5").Value = 19
$ws.Range("C5").Value = "60b45e9961dd412bfb6780f8"
$ws.Range("D5").Value = "Jewel"
$ws.Range("F5").Value = 8.390562821666926

# Row 6
$ws.Range("B6").Value = 21
$ws.Range("C6").Value = "5c0e89c6c323400001e6c4a5"
$ws.Range("D6").Value = "Bri"
$ws.Range("F6").Value = 8.284137808845447

# Row 7
$ws.Range("B7").Value = 32
$ws.Range("C7").Value = "6036f9b3b1842f8b659b18c7"
$ws.Range("D7").Value = "Kellie"
$ws.Range("F7").Value = 5.499920003737663
$ws.Range("G7").Value = "White"

# Row 8
$ws.Range("F8").Value = 5.496086788842061

# Row 9
$ws.Range("B9").Value = 33
$ws.Range("C9").Value = "60cb36ee9f58331a33cf5506"
$ws.Range("D9").Value = "Shaniek"
$ws.Range("F9").Value = 5.062422754775289
$ws.Range("G9").Value = "Black or African American"

# Row 10
$ws.Range("B10").Value = 35
$ws.Range("C10").Value = "6077db0613ce87b4a62a78f9"
$ws.Range("D10").Value = "Lori"
$ws.Range("F10").Value = 4.244814854093466

# Row 11
$ws.Range("B11").Value = 34
$ws.Range("C11").Value = "5e96194b0a9fe909389e9f7b"
$ws.Range("D11").Value = "Tina"
$ws.Range("F11").Value = 4.079969157910064

# Row 12
$ws.Range("F12").Value = 2.218831050136576

# Row 13
$ws.Range("F13").Value = 1.453411958882284
